$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the contact name in D3 from "Dakota Myers" to "Franz Ferdinand"
$ws.Range("D3").Value = "Franz Ferdinand"

# Move the active selection to the edited cell, matching the saved selection state
$ws.Range("D3").Select()
